$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials")

# --- Remove last row (#28, "Automation Scripts" / pypi:openpyxl entry no longer tracked here) ---
$ws.Rows.Item(29).Delete()

# --- Column widths: A=5, B=30, C=40, D=15, E=25, F=40 characters ---
# (ColumnWidth property is offset by ~0.8333 from the stored OOXML "width"; compensate so the
#  saved file shows the exact target widths.)
$ws.Columns.Item(1).ColumnWidth = 4.166666666666667
$ws.Columns.Item(2).ColumnWidth = 29.166666666666668
$ws.Columns.Item(3).ColumnWidth = 39.166666666666664
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 24.166666666666668
$ws.Columns.Item(6).ColumnWidth = 39.166666666666664

# --- Header row: F1 "Status" -> "Comment" ---
$ws.Range("F1").Value = "Comment"

# Row 2 (# 1): npm:react
$ws.Range("B2").Value = "Frontend (React UI)"
$ws.Range("C2").Value = "npm:react"
$ws.Range("D2").Value = "18.3.1"
$ws.Range("E2").Value = "MIT"
$ws.Range("F2").Value = ""

# Row 3 (# 2): npm:react-dom
$ws.Range("B3").Value = "Frontend (React UI)"
$ws.Range("C3").Value = "npm:react-dom"
$ws.Range("D3").Value = "18.3.1"
$ws.Range("E3").Value = "MIT"
$ws.Range("F3").Value = ""

# Row 4 (# 3): pypi:aiortc
$ws.Range("B4").Value = "Backend (FastAPI API)"
$ws.Range("C4").Value = "pypi:aiortc"
$ws.Range("D4").Value = "1.14.0"
$ws.Range("E4").Value = "BSD-3-Clause"
$ws.Range("F4").Value = ""

# Row 5 (# 4): pypi:av
$ws.Range("B5").Value = "Backend (FastAPI API)"
$ws.Range("C5").Value = "pypi:av"
$ws.Range("D5").Value = "16.0.1"
$ws.Range("E5").Value = "BSD-3-Clause"
$ws.Range("F5").Value = ""

# Row 6 (# 5): pypi:fastapi
$ws.Range("B6").Value = "Backend (FastAPI API)"
$ws.Range("C6").Value = "pypi:fastapi"
$ws.Range("D6").Value = "0.115.10"
$ws.Range("E6").Value = "MIT"
$ws.Range("F6").Value = ""

# Row 7 (# 6): pypi:httpx
$ws.Range("B7").Value = "Backend (FastAPI API)"
$ws.Range("C7").Value = "pypi:httpx"
$ws.Range("D7").Value = "0.27.2"
$ws.Range("E7").Value = "BSD License"
$ws.Range("F7").Value = ""

# Row 8 (# 7): pypi:numpy
$ws.Range("B8").Value = "Backend (FastAPI API)"
$ws.Range("C8").Value = "pypi:numpy"
$ws.Range("D8").Value = "1.26.4"
$ws.Range("E8").Value = "BSD-3-Clause"
$ws.Range("F8").Value = ""

# Row 9 (# 8): pypi:opencv-python
$ws.Range("B9").Value = "Backend (FastAPI API)"
$ws.Range("C9").Value = "pypi:opencv-python"
$ws.Range("D9").Value = "4.9.0.80"
$ws.Range("E9").Value = "Apache 2.0"
$ws.Range("F9").Value = ""

# Row 10 (# 9): pypi:pydantic
$ws.Range("B10").Value = "Backend (FastAPI API)"
$ws.Range("C10").Value = "pypi:pydantic"
$ws.Range("D10").Value = "2.12.3"
$ws.Range("E10").Value = "MIT"
$ws.Range("F10").Value = ""

# Row 11 (# 10): pypi:timm
$ws.Range("B11").Value = "Backend (FastAPI API)"
$ws.Range("C11").Value = "pypi:timm"
$ws.Range("D11").Value = "1.0.22"
$ws.Range("E11").Value = "Apache-2.0"
$ws.Range("F11").Value = ""

# Row 12 (# 11): pypi:ultralytics
$ws.Range("B12").Value = "Backend (FastAPI API)"
$ws.Range("C12").Value = "pypi:ultralytics"
$ws.Range("D12").Value = "'8.3.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "AGPL-3.0"
$ws.Range("F12").Value = ""

# Row 13 (# 12): pypi:uvicorn
$ws.Range("B13").Value = "Backend (FastAPI API)"
$ws.Range("C13").Value = "pypi:uvicorn"
$ws.Range("D13").Value = "0.38.0"
$ws.Range("E13").Value = "BSD-3-Clause"
$ws.Range("F13").Value = ""

# Row 14 (# 13): npm:@types/react
$ws.Range("B14").Value = "Frontend (React UI)"
$ws.Range("C14").Value = "npm:@types/react"
$ws.Range("D14").Value = "18.3.12"
$ws.Range("E14").Value = "MIT"
$ws.Range("F14").Value = ""

# Row 15 (# 14): npm:@types/react-dom
$ws.Range("B15").Value = "Frontend (React UI)"
$ws.Range("C15").Value = "npm:@types/react-dom"
$ws.Range("D15").Value = "18.3.1"
$ws.Range("E15").Value = "MIT"
$ws.Range("F15").Value = ""

# Row 16 (# 15): npm:@typescript-eslint/eslint-plugin
$ws.Range("B16").Value = "Frontend (React UI)"
$ws.Range("C16").Value = "npm:@typescript-eslint/eslint-plugin"
$ws.Range("D16").Value = "8.16.0"
$ws.Range("E16").Value = "MIT"
$ws.Range("F16").Value = ""

# Row 17 (# 16): npm:@typescript-eslint/parser
$ws.Range("B17").Value = "Frontend (React UI)"
$ws.Range("C17").Value = "npm:@typescript-eslint/parser"
$ws.Range("D17").Value = "8.16.0"
$ws.Range("E17").Value = "BSD-2-Clause"
$ws.Range("F17").Value = ""

# Row 18 (# 17): npm:@vitejs/plugin-react
$ws.Range("B18").Value = "Frontend (React UI)"
$ws.Range("C18").Value = "npm:@vitejs/plugin-react"
$ws.Range("D18").Value = "4.3.4"
$ws.Range("E18").Value = "MIT"
$ws.Range("F18").Value = ""

# Row 19 (# 18): npm:eslint
$ws.Range("B19").Value = "Frontend (React UI)"
$ws.Range("C19").Value = "npm:eslint"
$ws.Range("D19").Value = "9.17.0"
$ws.Range("E19").Value = "MIT"
$ws.Range("F19").Value = ""

# Row 20 (# 19): pypi:mypy
$ws.Range("B20").Value = "Backend Dev Dependencies"
$ws.Range("C20").Value = "pypi:mypy"
$ws.Range("D20").Value = "1.13.0"
$ws.Range("E20").Value = "MIT"
$ws.Range("F20").Value = ""

# Row 21 (# 20): npm:prettier
$ws.Range("B21").Value = "Frontend (React UI)"
$ws.Range("C21").Value = "npm:prettier"
$ws.Range("D21").Value = "3.4.2"
$ws.Range("E21").Value = "MIT"
$ws.Range("F21").Value = ""

# Row 22 (# 21): pypi:pytest
$ws.Range("B22").Value = "Backend Dev Dependencies"
$ws.Range("C22").Value = "pypi:pytest"
$ws.Range("D22").Value = "8.3.3"
$ws.Range("E22").Value = "MIT"
$ws.Range("F22").Value = ""

# Row 23 (# 22): pypi:pytest-asyncio
$ws.Range("B23").Value = "Backend Dev Dependencies"
$ws.Range("C23").Value = "pypi:pytest-asyncio"
$ws.Range("D23").Value = "1.2.0"
$ws.Range("E23").Value = "Apache-2.0"
$ws.Range("F23").Value = ""

# Row 24 (# 23): pypi:reuse
$ws.Range("B24").Value = "Backend Dev Dependencies"
$ws.Range("C24").Value = "pypi:reuse"
$ws.Range("D24").Value = "4.0.3"
$ws.Range("E24").Value = "Apache-2.0 AND CC0-1.0 AND CC-BY-SA-4.0 AND GPL-3.0-or-later"
$ws.Range("F24").Value = ""

# Row 25 (# 24): pypi:ruff
$ws.Range("B25").Value = "Backend Dev Dependencies"
$ws.Range("C25").Value = "pypi:ruff"
$ws.Range("D25").Value = "0.7.0"
$ws.Range("E25").Value = "MIT"
$ws.Range("F25").Value = ""

# Row 26 (# 25): npm:typescript
$ws.Range("B26").Value = "Frontend (React UI)"
$ws.Range("C26").Value = "npm:typescript"
$ws.Range("D26").Value = "5.6.3"
$ws.Range("E26").Value = "Apache-2.0"
$ws.Range("F26").Value = ""

# Row 27 (# 26): npm:vite
$ws.Range("B27").Value = "Frontend (React UI)"
$ws.Range("C27").Value = "npm:vite"
$ws.Range("D27").Value = "6.0.3"
$ws.Range("E27").Value = "MIT"
$ws.Range("F27").Value = ""

# Row 28 (# 27): npm:vitest
$ws.Range("B28").Value = "Frontend (React UI)"
$ws.Range("C28").Value = "npm:vitest"
$ws.Range("D28").Value = "2.1.5"
$ws.Range("E28").Value = "MIT"
$ws.Range("F28").Value = ""
